$d = $word.ActiveDocument

$replacements = @(
    @{old = "2024-08-11 Sunday"; new = "2024-08-12 Monday"},
    @{old = "44×99=4356"; new = "48×19=912"},
    @{old = "76×48=3648"; new = "64×24=1536"},
    @{old = "60×51=3060"; new = "20×36=720"},
    @{old = "94×56=5264"; new = "73×18=1314"},
    @{old = "52×42=2184"; new = "14×78=1092"},
    @{old = "38×85=3230"; new = "46×76=3496"},
    @{old = "93×42=3906"; new = "18×56=1008"},
    @{old = "15×26=390"; new = "66×43=2838"},
    @{old = "36×37=1332"; new = "46×96=4416"},
    @{old = "36×15=540"; new = "84×84=7056"},
    @{old = "77×30=2310"; new = "17×82=1394"},
    @{old = "64×98=6272"; new = "72×46=3312"},
    @{old = "83×73=6059"; new = "29×95=2755"},
    @{old = "41×11=451"; new = "77×98=7546"},
    @{old = "85×89=7565"; new = "76×12=912"},
    @{old = "81×64=5184"; new = "72×48=3456"},
    @{old = "57×60=3420"; new = "87×92=8004"},
    @{old = "35×32=1120"; new = "89×33=2937"},
    @{old = "89×80=7120"; new = "79×68=5372"},
    @{old = "26×33=858"; new = "53×42=2226"},
    @{old = "66×21=1386"; new = "21×63=1323"},
    @{old = "95×30=2850"; new = "46×40=1840"},
    @{old = "63×55=3465"; new = "78×27=2106"},
    @{old = "35×42=1470"; new = "30×27=810"},
    @{old = "61×58=3538"; new = "64×65=4160"}
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $pair.new, 2)
}

$d.Save()
